$d = $word.ActiveDocument

# The document title (paragraph 1, style "Titre") ends with the version number
# "V2.5". The edit bumps the minor version to "V2.6", and in doing so the final
# digit ends up in its own run (as if the "5" had been selected and retyped as
# "6"), leaving the rest of the title text in the original run.
$titlePara = $d.Paragraphs(1).Range

# Range.Text includes the trailing paragraph mark, so trim it off before doing
# any text comparisons / length math.
$rawText = $titlePara.Text
$titleText = $rawText.TrimEnd([char]13, [char]7)

$oldSuffix = "V2.5"
$newChar = "6"

if (-not $titleText.EndsWith($oldSuffix)) {
  throw "Unexpected title text, cannot locate '$oldSuffix': $titleText"
}

$paraStart = $titlePara.Start
$textEnd = $paraStart + $titleText.Length   # end of real text, excludes paragraph mark
$splitPos = $textEnd - 1                    # position right before the final digit ("5")

# Recover the formatting attributes (e.g. w:rsidRPr) carried by the original run
# so the retained portion of the text keeps them unchanged.
$openXml = $titlePara.WordOpenXML
$runAttrs = ""
if ($openXml -match '<w:r\s*([^>]*)><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>') {
  $runAttrs = $matches[1].Trim()
}
$attrXml = ""
if ($runAttrs.Length -gt 0) {
  $attrXml = " " + $runAttrs
}

$beforeRange = $d.Range($paraStart, $splitPos)
$beforeText = $beforeRange.Text
$beforeText = $beforeText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

$pkgTemplate = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r{0}><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>{1}</w:t></w:r><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>{2}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$xml = $pkgTemplate -f $attrXml, $beforeText, $newChar

$fullRunRange = $d.Range($paraStart, $textEnd)
$fullRunRange.InsertXML($xml)
